$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44211
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 6000
$ws.Range('R2').Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 3000

# Row 3
$ws.Range("D3").Value = 44592
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range('R3').Value = 'Provincia de Linares'
$ws.Range("S3").Value = 4000

# Row 5
$ws.Range("D5").Value = 44585
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 6500
$ws.Range("P5").Value = 6750
$ws.Range("S5").Value = 3375

# Row 6
$ws.Range("D6").Value = 44214
$ws.Range("M6").Value = 48
$ws.Range('R6').Value = 'Provincia de Linares'

# Row 7
$ws.Range("D7").Value = 44589
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range('R7').Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 3000

# Row 8
$ws.Range("D8").Value = 44628
$ws.Range("M8").Value = 40

# Row 9
$ws.Range("D9").Value = 44582
$ws.Range("M9").Value = 150
$ws.Range("O9").Value = 6500
$ws.Range("P9").Value = 6233
$ws.Range('R9').Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 3116

# Row 10
$ws.Range("D10").Value = 44587
$ws.Range("M10").Value = 165
$ws.Range("N10").Value = 6500
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6742
$ws.Range('R10').Value = 'Provincia de Linares'
$ws.Range("S10").Value = 3371

# Row 11
$ws.Range("D11").Value = 44614
$ws.Range("M11").Value = 45
$ws.Range('R11').Value = 'Provincia de Linares'

# Row 12
$ws.Range("D12").Value = 44588

# Row 13
$ws.Range("D13").Value = 44959
$ws.Range("N13").Value = 7000
$ws.Range("O13").Value = 7000
$ws.Range("P13").Value = 7000
$ws.Range('R13').Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 3500

# Row 14
$ws.Range("D14").Value = 44974
$ws.Range("M14").Value = 130
$ws.Range("O14").Value = 7500
$ws.Range("P14").Value = 7269
$ws.Range("S14").Value = 3634

# Row 15
$ws.Range("D15").Value = 44586
$ws.Range("M15").Value = 80

# Row 16
$ws.Range("D16").Value = 44606
$ws.Range("M16").Value = 45
$ws.Range("N16").Value = 7000
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 7000
$ws.Range("S16").Value = 3500

# Row 17
$ws.Range("D17").Value = 44960
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 7000
$ws.Range("P17").Value = 7000
$ws.Range("S17").Value = 3500

# Row 18
$ws.Range("D18").Value = 44209
$ws.Range("M18").Value = 58
$ws.Range("N18").Value = 6000
$ws.Range("O18").Value = 6000
$ws.Range("P18").Value = 6000
$ws.Range("S18").Value = 3000

